$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the branch name text in A3 (shared string text edit)
$ws.Range("A3").Value = "CIMS.CAN.AT.Residential.Dwellings.Lighting"

# Replace the shared-formula chain in M3:W3 with literal values
$ws.Range("M3").Value = 1
$ws.Range("N3").Value = 1
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 0.5
$ws.Range("Q3").Value = 0.0001
$ws.Range("R3").Value = 0.0001
$ws.Range("S3").Value = 0.0001
$ws.Range("T3").Value = 0.0001
$ws.Range("U3").Value = 0.0001
$ws.Range("V3").Value = 0.0001
$ws.Range("W3").Value = 0.0001

# Update the sheet selection to match the author's saved view state
$ws.Range("A1:X4").Select() | Out-Null
